$d = $word.ActiveDocument

# Helper-less approach: we manually track the insertion position ($pos) and,
# for each fragment, insert the text right after that position using a
# collapsed Range + InsertAfter, then (optionally) apply a character style
# to the exact Range that was just inserted. This produces separate <w:r>
# runs with distinct <w:rStyle> values, matching the target OOXML.

function Insert-StyledText {
    param(
        [int]$Position,
        [string]$Text,
        [string]$StyleName
    )
    $ip = $d.Range($Position, $Position)
    $ip.InsertAfter($Text)
    if ($StyleName) {
        $r = $d.Range($Position, $Position + $Text.Length)
        $r.Style = $StyleName
    }
    return $Position + $Text.Length
}

function Insert-LineBreak {
    param(
        [int]$Position
    )
    $ip = $d.Range($Position, $Position)
    $ip.InsertBreak(6)  # wdLineBreak
    return $Position + 1
}

# NOTE: this runtime's PowerShell subset does not bind named parameters
# (e.g. `-Position $x`) to `param()` blocks, so all calls below pass
# arguments positionally.

# -----------------------------------------------------------------------
# Hunk 1: "nchar(chr)" -> "nchar(chr); stringr::str_count(chr)"
# Insert "(chr); stringr" + "::" + "str_count" right before the existing
# "(chr)" run that follows "nchar".
# -----------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("nchar(chr)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos = $rng.Start + 5   # length of "nchar" -> position right before "(chr)"

$pos = Insert-StyledText $pos "(chr); stringr" "NormalTok"
$pos = Insert-StyledText $pos "::" "OperatorTok"
$pos = Insert-StyledText $pos "str_count" "KeywordTok"

# -----------------------------------------------------------------------
# Hunk 2: after "...paste0(chr_a, chr_b, collapse = ...) # with no separator"
# (the second such comment, right before "# split chr_a ...") insert a new
# comment line plus a stringr::str_c(...) example line.
# -----------------------------------------------------------------------

$rng2 = $d.Content
$rng2.Find.Execute("collapse = ...) # with no separator", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$pos2 = $rng2.End

$pos2 = Insert-LineBreak $pos2
$pos2 = Insert-StyledText $pos2 "# same functionality in stringr" "CommentTok"
$pos2 = Insert-LineBreak $pos2
$pos2 = Insert-StyledText $pos2 "stringr" "NormalTok"
$pos2 = Insert-StyledText $pos2 "::" "OperatorTok"
$pos2 = Insert-StyledText $pos2 "str_c" "KeywordTok"
$pos2 = Insert-StyledText $pos2 "(chr_a, chr_b, " "NormalTok"
$pos2 = Insert-StyledText $pos2 "sep =" "DataTypeTok"
$pos2 = Insert-StyledText $pos2 " ..., " "NormalTok"
$pos2 = Insert-StyledText $pos2 "collapse =" "DataTypeTok"
$pos2 = Insert-StyledText $pos2 " ...)" "NormalTok"

Write-Host "Edits applied successfully"
